# Bætti við tímum í dagbók
# Adds new time entries (hours) for Sunday (column I) in the "Vika 4"
# (week 4) block of the time log, rows 16-21, and moves the active
# selection to I22. The dependent SUM formulas in column J (per-row
# totals), J22 (week total) and the summary table (D47:D53) recalculate
# automatically because they are formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rannsóknir (row 16) - sun
$ws.Range("I16").Value = 1
# Kröfulýsing (row 17) - sun
$ws.Range("I17").Value = 1
# Hönnun (row 18) - sun
$ws.Range("I18").Value = 2
# Frágangur (row 21) - sun
$ws.Range("I21").Value = 1.25

# Move the selected/active cell, matching the author's final click position
$ws.Range("I22").Select()
